# Testing with real data: rename/relabel the impedance & frequency headers.
# Column data (A:G) itself stays where it is; only the header text for the
# "z_real", "z_imag" and "frequency" columns changes (capitalization + wording),
# and the active selection moves to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set C1 first so the new shared-string "angular frequency" is appended
# before "Z_real"/"Z_imag", matching the authored workbook's string order.
$ws.Range("C1").Value = "angular frequency"
$ws.Range("A1").Value = "Z_real"
$ws.Range("B1").Value = "Z_imag"

# D1/E1/F1/G1 (eff_cap, applied voltage, J_ph, J) are unchanged.

# Update the active selection to B1, as in the saved workbook.
$ws.Range("B1").Select()
